$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price (D) column cells to Text format before writing so
# numeric-looking strings (e.g. "1.00", "0.0794", "245.00") are stored
# verbatim as text instead of being auto-coerced into numbers (which
# would drop trailing zeros / use floating point representations).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values scraped for this run.
$ws.Range("D2").Value = "45.889.11"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "2.357.51"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "301.63"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "99.43"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("D10").Value = "34.90"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("D12").Value = "7.17"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "2.711.62"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "2.378.23"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "13.78"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "0.811"
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("D18").Value = "45.857.54"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -7.20%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").Value = "66.46"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("D23").Value = "245.00"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").Value = "2.82"
$ws.Range("E24").Value = "  -5.64%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "1.91"
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").Value = "40.77"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "21.04"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").Value = "3.66"
$ws.Range("E31").Value = "  +16.17%  "
$ws.Range("D32").Value = "2.75"
$ws.Range("E32").Value = "  +5.20%  "
$ws.Range("E33").Value = "  -6.87%  "
$ws.Range("D34").Value = "145.19"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "15.23"
$ws.Range("E39").Value = "  +8.33%  "
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("D42").Value = "3.21"
$ws.Range("E42").Value = "  -5.60%  "
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "1.848.89"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "90.93"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "1.80"
$ws.Range("E46").Value = "  -8.17%  "
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("D48").Value = "70.68"
$ws.Range("E48").Value = "  -5.94%  "
$ws.Range("D49").Value = "2.583.17"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "8.04"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "96.31"
$ws.Range("E51").Value = "  -2.79%  "
